# Update OPP stats for Los Angeles Lakers (C), Los Angeles Clippers (E),
# Utah Jazz (I) and New Orleans Pelicans (S) columns with July 20 game data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 32
$ws.Range("E2").Value = 32
$ws.Range("I2").Value = 40
$ws.Range("S2").Value = 37

$ws.Range("C3").Value = 74
$ws.Range("E3").Value = 82
$ws.Range("I3").Value = 94
$ws.Range("S3").Value = 84

$ws.Range("C4").Value = 43.2
$ws.Range("E4").Value = 39
$ws.Range("I4").Value = 42.6
$ws.Range("S4").Value = 44

$ws.Range("C5").Value = 16
$ws.Range("E5").Value = 11
$ws.Range("I5").Value = 11
$ws.Range("S5").Value = 8

$ws.Range("C6").Value = 36
$ws.Range("E6").Value = 36
$ws.Range("I6").Value = 31
$ws.Range("S6").Value = 34

$ws.Range("C7").Value = 44.4
$ws.Range("E7").Value = 30.6
$ws.Range("I7").Value = 35.5
$ws.Range("S7").Value = 23.5

$ws.Range("C8").Value = 21
$ws.Range("E8").Value = 28
$ws.Range("I8").Value = 13
$ws.Range("S8").Value = 24

$ws.Range("C9").Value = 28
$ws.Range("E9").Value = 37
$ws.Range("I9").Value = 18
$ws.Range("S9").Value = 28

$ws.Range("C10").Value = 75
$ws.Range("E10").Value = 75.7
$ws.Range("I10").Value = 72.2
$ws.Range("S10").Value = 85.7

$ws.Range("C11").Value = 4
$ws.Range("E11").Value = 11
$ws.Range("I11").Value = 15
$ws.Range("S11").Value = 12

$ws.Range("C12").Value = 32
$ws.Range("E12").Value = 34
$ws.Range("I12").Value = 30
$ws.Range("S12").Value = 31

$ws.Range("C13").Value = 36
$ws.Range("E13").Value = 45
$ws.Range("I13").Value = 45
$ws.Range("S13").Value = 43

$ws.Range("C14").Value = 17
$ws.Range("E14").Value = 21
$ws.Range("I14").Value = 22
$ws.Range("S14").Value = 17

$ws.Range("C15").Value = 22
$ws.Range("E15").Value = 16
$ws.Range("I15").Value = 21
$ws.Range("S15").Value = 20

$ws.Range("C16").Value = 6
$ws.Range("E16").Value = 6
$ws.Range("I16").Value = 12
$ws.Range("S16").Value = 11

$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 3
$ws.Range("I17").Value = 2
$ws.Range("S17").Value = 6

$ws.Range("C18").Value = 3
$ws.Range("I18").Value = 6
$ws.Range("S18").Value = 2

$ws.Range("C19").Value = 30
$ws.Range("E19").Value = 27
$ws.Range("I19").Value = 25
$ws.Range("S19").Value = 23

$ws.Range("C20").Value = 27
$ws.Range("E20").Value = 30
$ws.Range("I20").Value = 23
$ws.Range("S20").Value = 25

$ws.Range("C21").Value = 101
$ws.Range("E21").Value = 103
$ws.Range("I21").Value = 104
$ws.Range("S21").Value = 106

# Trim the blank placeholder rows that are no longer needed at the bottom of
# the sheet, without shifting the rows that are being kept.
$ws.Range("A25:O41").Clear() | Out-Null
$ws.Range("A46:O49").Clear() | Out-Null
$ws.Range("A42:D44").Clear() | Out-Null

# Move the lone trailing cell from row 51 up to row 47.
$ws.Range("A51").Copy($ws.Range("A47")) | Out-Null
$ws.Range("A51").Clear() | Out-Null

# Restore the selection the author left the sheet with.
$ws.Range("A25:XFD28").Select() | Out-Null
